# Fixed bug where prev_row in process_trades() was not being reset between tests.
# Populate the TestCases sheet (Sheet1) with 6 test-case rows (rows 2-7),
# replacing the single existing test row + blank template rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 3-5 get the full A:K formatting (incl. the optional-settings column K,
# cleared of its value) copied down from row 2; rows 6-7 only need A:J.
$ws.Range("A2:K2").Copy($ws.Range("A3:K3"))
$ws.Range("A2:K2").Copy($ws.Range("A4:K4"))
$ws.Range("A2:K2").Copy($ws.Range("A5:K5"))
$ws.Range("A2:J2").Copy($ws.Range("A6:J6"))
$ws.Range("A2:J2").Copy($ws.Range("A7:J7"))

$exchange = "Bybit"
$pair     = "BTCUSDT"
$strategy = "MACD"
$from     = 44197
$to       = 44561

# A, interval(F), TP%(G), SL%(H), exit strategy(J)
$rows = @(
  @(1, "15m", 1.2, 1.2, "FixedPCT"),
  @(2, "15m", 1.5, 1.5, "ExitOnNextEntry"),
  @(3, "30m", 1.2, 1.2, "FixedPCT"),
  @(4, "30m", 1.5, 1.5, "ExitOnNextEntry"),
  @(5, "1h",  1.2, 1.3, "FixedPCT"),
  @(6, "1h",  1.5, 1.5, "ExitOnNextEntry")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $exchange
    $ws.Cells.Item($r, 3).Value = $pair
    $ws.Cells.Item($r, 4).Value = $from
    $ws.Cells.Item($r, 5).Value = $to
    $ws.Cells.Item($r, 6).Value = $row[1]
    $ws.Cells.Item($r, 7).Value = $row[2]
    $ws.Cells.Item($r, 8).Value = $row[3]
    $ws.Cells.Item($r, 9).Value = $strategy
    $ws.Cells.Item($r, 10).Value = $row[4]
    $r++
}

# The optional strategy-settings value only remains (blank, formatted) on
# rows 2-5; clear the text that used to live in K2.
$ws.Range("K2:K5").ClearContents()

# Update the active selection to H7, matching the saved workbook state.
$ws.Range("H7").Select()
